# Daily attendance processing - 2025-11-24 06:36:16
#
# Applies the day's refreshed attendance-report data to the "Session
# Analysis Results" sheet:
#   - "Recorded By" address lists re-ordered (same people, new order as
#     exported by the upstream system on this run)
#   - Session 1 of MICROBIOLOGY (row 12) flips from "Pending" to
#     "Not Recorded" (its date has passed with nothing recorded), so it
#     picks up the same pink "Not Recorded" look already used elsewhere
#     on the sheet
#   - The dependent Missing/Pending counters (both the quick-stats block
#     and the per-group summary row) shift by one to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded-by email lists: reordered by the source system ---------
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# --- Row 12 (MICROBIOLOGY, session 1): Pending -> Not Recorded -------
# Pick up the same fill/font formatting already used for "Not Recorded"
# rows (e.g. row 29) instead of inventing a fresh style.
$ws.Range("A29:I29").Copy() | Out-Null
$ws.Range("A12:I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I12").Value = "Not Recorded"

# --- Dependent counters: one more missing session, one less pending --
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 16
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 16
